{"js": "// Replace the name \"Luz\" with \"Bella\" in the greeting\n// \"Hi, Luz. How are you? \" -> \"Hi, Bella. How are you? \"\nconst body = context.document.body;\n\nconst results = body.search(\"Luz\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items,text\");\nawait context.sync();\n\nfor (const range of results.items) {\n  range.insertText(\"Bella\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the name \"Luz\" with \"Bella\" in the greeting:\n# \"Hi, Luz. How are you? \" -> \"Hi, Bella. How are you? \"\n$d = $word.ActiveDocument\n\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards,\n#              MatchSoundsLike, MatchAllWordForms, Forward, Wrap,\n#              Format, ReplaceWith, Replace)\n# Wrap = 1 (wdFindContinue), Replace = 2 (wdReplaceAll)\n$d.Content.Find.Execute(\"Luz\", $true, $true, $false, $false, $false, $true, 1, $false, \"Bella\", 2) | Out-Null\n"}
